$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Add a new row for a 4th product ("Bonelo", dog food) with its price/weight/stock.
$ws.Range("A5").Value = 4
$ws.Range("B5").Value = "Bonelo"
$ws.Range("C5").Value = "perro"
$ws.Range("D5").Value = 10600
$ws.Range("E5").Value = 20
$ws.Range("F5").Value = 10
$ws.Rows.Item(5).RowHeight = 15.75

# Rename the "alimento"/"precio" headers to their English equivalents.
$ws.Range("B1").Value = "name"
$ws.Range("D1").Value = "price"

# Leave the selection where the editor ended up.
$ws.Range("H15").Select()
